# FAXINAL_DO_SOTURNO.xlsx update
#   - "Paineis DARQ"              -> "PAINEIS DARQ"
#   - "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   - remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Rename sheets to the new uppercase/accented titles
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely (suppress the
# "delete sheet" confirmation prompt Excel would otherwise raise)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Keep the original active tab (first sheet) selected, since deleting the
# last sheet would otherwise shift the active tab onto its neighbour
$wb.Worksheets.Item("PAINEIS DARQ").Activate() | Out-Null
